# Atualizacao de bases das ligas - corrige a ordem de 3 jogos (linhas 88-90)
# e troca 2 jogos (linhas 118-119) que haviam sido gravados na ordem errada.
# Mantem a coluna A (indice sequencial do jogo) e a coluna D (data) intactas -
# apenas os dados de B:AD (id do jogo, times, placares, odds, etc.) trocam de
# linha.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Linhas 88, 89, 90: rotacao ciclica -----------------------------------
# novo(88) = antigo(90); novo(89) = antigo(88); novo(90) = antigo(89)
$row88 = $ws.Range("B88:AD88").Value()
$row89 = $ws.Range("B89:AD89").Value()
$row90 = $ws.Range("B90:AD90").Value()

$ws.Range("B88:AD88").Value = $row90
$ws.Range("B89:AD89").Value = $row88
$ws.Range("B90:AD90").Value = $row89

# --- Linhas 118, 119: troca simples ---------------------------------------
$row118 = $ws.Range("B118:AD118").Value()
$row119 = $ws.Range("B119:AD119").Value()

$ws.Range("B118:AD118").Value = $row119
$ws.Range("B119:AD119").Value = $row118
